# Update the speedup columns (G:I) on the active sheet so that they are
# computed as "no optimization" (column B) divided by each optimization
# stage (columns C, D, E) instead of the previous "baseline" (column F,
# which is always 1) divided by the optimization stage columns.
#
# G = B / C   (VP_speedup)
# H = B / D   (VP+EF_speedup)
# I = B / E   (VP+EF+RL_speedup)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2  # column B - no optimization
    $c = $ws.Cells.Item($r, 3).Value2  # column C - VP
    $d = $ws.Cells.Item($r, 4).Value2  # column D - VP+EF
    $e = $ws.Cells.Item($r, 5).Value2  # column E - VP+EF+RL

    $ws.Cells.Item($r, 7).Value = $b / $c   # column G - VP_speedup
    $ws.Cells.Item($r, 8).Value = $b / $d   # column H - VP+EF_speedup
    $ws.Cells.Item($r, 9).Value = $b / $e   # column I - VP+EF+RL_speedup
}
